$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New text entries added to the sheet (documentation / ideas).
# Entered in the same order the author typed them so the shared-string
# table ends up in the same order as the source edit.
$ws.Range("F13").Value = "Back to the drawing board"
$ws.Range("C15").Value = "It burns"
$ws.Range("C16").Value = "A race with death"
$ws.Range("F14").Value = "And then it stopped"
$ws.Range("H5").Value = "For death and glory"

# Update selection to match the final state of the diff
$ws.Range("H16").Select()
